# Translate the ContosoLearn Market Research document body from
# English to Italian, paragraph by paragraph. Each paragraph's whole
# text (its Range, excluding the trailing paragraph mark) is replaced
# in one shot via Range.Text assignment -- this is the COM-interop
# equivalent of retyping the paragraph, keeps straight apostrophes
# (no autocorrect smart-quote substitution), and naturally collapses
# into a single run when (as here) run formatting is uniform, which
# also has the effect of dropping the stray <w:proofErr/> markers that
# used to split the WoodgroveLMS sentence into three runs.

$d = $word.ActiveDocument

function Set-ParagraphText($index, $text) {
    $p = $d.Paragraphs.Item($index)
    $r = $d.Range($p.Range.Start, $p.Range.End - 1)
    $r.Text = $text
}

Set-ParagraphText 1 "Ricerca di mercato di ContosoLearn"

$adatumText = "AdatumLearn: AdatumLearn è una delle principali piattaforme di apprendimento basate sull'intelligenza artificiale, che arricchisce l'eLearning automatizzando un'ampia gamma di attività. " + "È nota per le sue funzionalità di creazione di contenuti e per la tecnologia di apprendimento adattivo."
Set-ParagraphText 2 $adatumText

Set-ParagraphText 3 "AdventureLearn: AdventureLearn è un'altra piattaforma di apprendimento basata sull'intelligenza artificiale, che offre esperienze di apprendimento personalizzate e raccomandazioni basate sull'analisi dei dati."

Set-ParagraphText 4 "AlpineTraining: AlpineTraining è una piattaforma di apprendimento ottimizzata per l'uso su dispositivi mobili, focalizzata sul microlearning."

Set-ParagraphText 5 "Bellows OnDemand: Bellows OnDemand offre una soluzione di apprendimento completa, comprendente strumenti per la creazione di contenuti e la collaborazione sociale."

Set-ParagraphText 6 "FabrikamLearning: FabrikamLearning fornisce una suite di piattaforme in grado di rispondere a diverse esigenze di apprendimento."

Set-ParagraphText 7 "FirstUp Cards: FirstUp Cards è un'app di apprendimento mobile ideale per gestire training su procedure di sicurezza, conformità, conoscenza di nuovi prodotti o qualsiasi altra tipologia di scenario di apprendimento."

Set-ParagraphText 8 "Munson'sLearn: Munson'sLearn supporta le aziende nella formazione di dipendenti, partner e clienti."

Set-ParagraphText 9 "LibertyLearn: LibertyLearn è un LMS veloce e progettato per progetti mission-critical."

Set-ParagraphText 10 "WoodgroveLMS: WoodgroveLMS è un sistema di gestione dell'apprendimento funzionale ed esteticamente curato, concepito per offrire un'esperienza didattica di alto livello."

Set-ParagraphText 11 "NorthwindWorlds: NorthwindWorlds è una soluzione di apprendimento potente, intuitiva e affidabile, adatta sia a singoli individui che a realtà aziendali."

Set-ParagraphText 12 "ProsewareLearn: ProsewareLearn è un'azienda di formazione online che propone un'ampia selezione di corsi video per sviluppatori software, amministratori IT e professionisti creativi."

Set-ParagraphText 13 "RelecloudLearn: RelecloudLearn è una piattaforma di apprendimento online che offre MOOC (corsi aperti su larga scala), specializzazioni e titoli di studio in numerose discipline."

Set-ParagraphText 14 "TreyAcademy: TreyAcademy è una piattaforma di apprendimento online, lanciata nel maggio 2010, rivolta a professionisti e studenti."

$finalText = "Queste piattaforme vantano una notevole presenza sul mercato e sono ampiamente riconosciute per le loro funzionalità basate sull'AI, tra cui esperienze di apprendimento personalizzate, raccomandazioni guidate dai dati e automazione dei processi. " + "Stanno trasformando il panorama di eLearning, sfruttando l'intelligenza artificiale per offrire esperienze di apprendimento più coinvolgenti, gratificanti e su misura per ogni utente. "
Set-ParagraphText 15 $finalText
